# Add new flag definitions for "s" (suspect behavior) and "e" (data errant) under the flag_met section.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing the flag_met block (currently starting at row 18) down.
$ws.Rows.Item(16).Insert()

$ws.Cells.Item(16, 1).Value = "s"
$ws.Cells.Item(16, 2).Value = "suspect behavior, use with caution"
$ws.Cells.Item(17, 1).Value = "e"
$ws.Cells.Item(17, 2).Value = "data errant, recoded to NA"

# Update the selected cell to reflect the edit location.
$ws.Range("D17").Select()
